# Apply scheduled-runner price/profit updates to Sheets (per commit diff)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 9
$ws.Range("H9").Value = 103.77778
$ws.Range("I9").Value = 105.1
$ws.Range("K9").Value = 105.1
$ws.Range("M9").Value = 63.90000000000001

# ALC row 15
$ws.Range("H15").Value = 1081.04
$ws.Range("I15").Value = 1081.04
$ws.Range("K15").Value = 3243.12
$ws.Range("M15").Value = -3074.12

# ALC row 107
$ws.Range("H107").Value = 2038.5
$ws.Range("J107").Value = 1918.7142
$ws.Range("L107").Value = 1918.7142
$ws.Range("N107").Value = -5758.7142

# ALC row 113
$ws.Range("H113").Value = 5527.7104
$ws.Range("J113").Value = 5163.3066
$ws.Range("L113").Value = 5163.3066
$ws.Range("N113").Value = -11671.3066

# ALC row 131
$ws.Range("H131").Value = 2439.0833
$ws.Range("I131").Value = 2542.6365
$ws.Range("K131").Value = 7627.9095
$ws.Range("M131").Value = -2587.9095

# ALC row 132
$ws.Range("H132").Value = 1912.4445
$ws.Range("I132").Value = 1774.8049
$ws.Range("J132").Value = 3323.25
$ws.Range("K132").Value = 5324.4147
$ws.Range("L132").Value = 9969.75
$ws.Range("M132").Value = -2794.4147
$ws.Range("N132").Value = -15029.75

# ALC row 137
$ws.Range("H137").Value = 2575.2
$ws.Range("I137").Value = 3311.9092
$ws.Range("J137").Value = 1674.7778
$ws.Range("K137").Value = 9935.7276
$ws.Range("L137").Value = 5024.3334
$ws.Range("M137").Value = -7385.7276
$ws.Range("N137").Value = -10124.3334

# ALC row 138
$ws.Range("H138").Value = 18870232
$ws.Range("J138").Value = 2654
$ws.Range("L138").Value = 7962
$ws.Range("N138").Value = -18242

$ws = $wb.Worksheets.Item("ARM")
# ARM row 32
$ws.Range("H32").Value = 2801.612
$ws.Range("I32").Value = 1716.3934
$ws.Range("K32").Value = 1716.3934
$ws.Range("M32").Value = -1429.3934

# ARM row 74
$ws.Range("H74").Value = 2887.375
$ws.Range("I74").Value = 2289.24
$ws.Range("J74").Value = 3884.2666
$ws.Range("K74").Value = 2289.24
$ws.Range("L74").Value = 3884.2666
$ws.Range("M74").Value = -1415.24
$ws.Range("N74").Value = -5632.2666

# ARM row 77
$ws.Range("H77").Value = 2887.375
$ws.Range("I77").Value = 2289.24
$ws.Range("J77").Value = 3884.2666
$ws.Range("K77").Value = 11446.2
$ws.Range("L77").Value = 19421.333
$ws.Range("M77").Value = -7078.199999999999
$ws.Range("N77").Value = -28157.333

# ARM row 88
$ws.Range("H88").Value = 94271.91
$ws.Range("I88").Value = 251249.75
$ws.Range("J88").Value = 4570.2856
$ws.Range("K88").Value = 251249.75
$ws.Range("L88").Value = 4570.2856
$ws.Range("M88").Value = -250843.75
$ws.Range("N88").Value = -5382.2856

# ARM row 91
$ws.Range("H91").Value = 94271.91
$ws.Range("I91").Value = 251249.75
$ws.Range("J91").Value = 4570.2856
$ws.Range("K91").Value = 251249.75
$ws.Range("L91").Value = 4570.2856
$ws.Range("M91").Value = -249845.75
$ws.Range("N91").Value = -7378.2856

# ARM row 122
$ws.Range("H122").Value = 5335.6665
$ws.Range("I122").Value = 5351.1816
$ws.Range("J122").Value = 5293
$ws.Range("K122").Value = 16053.5448
$ws.Range("L122").Value = 15879
$ws.Range("M122").Value = -13603.5448
$ws.Range("N122").Value = -20779

# ARM row 132
$ws.Range("H132").Value = 31261.07
$ws.Range("I132").Value = 32470.31
$ws.Range("J132").Value = 23126.182
$ws.Range("K132").Value = 97410.93000000001
$ws.Range("L132").Value = 69378.546
$ws.Range("M132").Value = -94880.93000000001
$ws.Range("N132").Value = -74438.546

$ws = $wb.Worksheets.Item("BSM")
# BSM row 20
$ws.Range("H20").Value = 4015.8
$ws.Range("I20").Value = 4399.6665
$ws.Range("J20").Value = 3440
$ws.Range("K20").Value = 4399.6665
$ws.Range("L20").Value = 3440
$ws.Range("M20").Value = -4152.6665
$ws.Range("N20").Value = -3934

# BSM row 86
$ws.Range("H86").Value = 3124.9285
$ws.Range("I86").Value = 2401.125
$ws.Range("K86").Value = 2401.125
$ws.Range("M86").Value = -1278.125

# BSM row 89
$ws.Range("H89").Value = 3124.9285
$ws.Range("I89").Value = 2401.125
$ws.Range("K89").Value = 12005.625
$ws.Range("M89").Value = -6389.625

# BSM row 128
$ws.Range("H128").Value = 2970
$ws.Range("I128").Value = 2970
$ws.Range("K128").Value = 8910
$ws.Range("M128").Value = -6420

$ws = $wb.Worksheets.Item("CRP")
# CRP row 31
$ws.Range("H31").Value = 1786.6
$ws.Range("I31").Value = 1570.3334
$ws.Range("J31").Value = 2111
$ws.Range("K31").Value = 1570.3334
$ws.Range("L31").Value = 2111
$ws.Range("M31").Value = -1275.3334
$ws.Range("N31").Value = -2701

# CRP row 34
$ws.Range("H34").Value = 1786.6
$ws.Range("I34").Value = 1570.3334
$ws.Range("J34").Value = 2111
$ws.Range("K34").Value = 1570.3334
$ws.Range("L34").Value = 2111
$ws.Range("M34").Value = -1368.3334
$ws.Range("N34").Value = -2515

# CRP row 58
$ws.Range("H58").Value = 2467.2144
$ws.Range("I58").Value = 2204.2
$ws.Range("K58").Value = 2204.2
$ws.Range("M58").Value = -2001.2

# CRP row 136
$ws.Range("H136").Value = 2467.2144
$ws.Range("I136").Value = 2204.2
$ws.Range("K136").Value = 6612.599999999999
$ws.Range("M136").Value = -4062.599999999999

$ws = $wb.Worksheets.Item("CUL")
# CUL row 9
$ws.Range("H9").Value = 5350
$ws.Range("I9").Value = 2200
$ws.Range("K9").Value = 6600
$ws.Range("M9").Value = -6376

# CUL row 49
$ws.Range("H49").Value = 1000000
$ws.Range("I49").Value = 1000000
$ws.Range("J49").Value = 0
$ws.Range("K49").Value = 3000000
$ws.Range("L49").Value = 0
$ws.Range("M49").Value = -2999844
$ws.Range("N49").ClearContents()

# CUL row 55
$ws.Range("H55").Value = 9455.556
$ws.Range("I55").Value = 2020
$ws.Range("J55").Value = 18750
$ws.Range("K55").Value = 6060
$ws.Range("L55").Value = 56250
$ws.Range("M55").Value = -5883
$ws.Range("N55").Value = -56604

# CUL row 116
$ws.Range("H116").Value = 1842.6666
$ws.Range("I116").Value = 1842.6666
$ws.Range("K116").Value = 5527.9998
$ws.Range("M116").Value = -2085.9998

# CUL row 117
$ws.Range("H117").Value = 88933.086
$ws.Range("J117").Value = 96981.45
$ws.Range("L117").Value = 290944.35
$ws.Range("N117").Value = -297828.35

# CUL row 119
$ws.Range("H119").Value = 999.5
$ws.Range("I119").Value = 999.5
$ws.Range("K119").Value = 2998.5
$ws.Range("M119").Value = 1839.5

# CUL row 120
$ws.Range("H120").Value = 14400
$ws.Range("J120").Value = 16000
$ws.Range("L120").Value = 48000
$ws.Range("N120").Value = -57676

# CUL row 123
$ws.Range("H123").Value = 6880
$ws.Range("I123").Value = 640
$ws.Range("K123").Value = 1920
$ws.Range("M123").Value = 530

# CUL row 131
$ws.Range("H131").Value = 1031.6
$ws.Range("I131").Value = 666
$ws.Range("J131").Value = 2494
$ws.Range("K131").Value = 1998
$ws.Range("L131").Value = 7482
$ws.Range("M131").Value = 3042
$ws.Range("N131").Value = -17562

$ws = $wb.Worksheets.Item("GSM")
# GSM row 80
$ws.Range("H80").Value = 3574.4614
$ws.Range("I80").Value = 3340.7144
$ws.Range("K80").Value = 3340.7144
$ws.Range("M80").Value = -2342.7144

# GSM row 83
$ws.Range("H83").Value = 3574.4614
$ws.Range("I83").Value = 3340.7144
$ws.Range("K83").Value = 16703.572
$ws.Range("M83").Value = -11711.572

# GSM row 132
$ws.Range("H132").Value = 1765.8334
$ws.Range("I132").Value = 1765.8334
$ws.Range("K132").Value = 5297.5002
$ws.Range("M132").Value = -2767.5002

$ws = $wb.Worksheets.Item("LTW")
# LTW row 39
$ws.Range("H39").Value = 5000
$ws.Range("J39").Value = 5000
$ws.Range("L39").Value = 5000
$ws.Range("N39").Value = -5920

# LTW row 122
$ws.Range("H122").Value = 5487.136
$ws.Range("I122").Value = 4192.5386
$ws.Range("K122").Value = 12577.6158
$ws.Range("M122").Value = -10127.6158

# LTW row 132
$ws.Range("H132").Value = 3623.0435
$ws.Range("I132").Value = 3154.0667
$ws.Range("K132").Value = 9462.2001
$ws.Range("M132").Value = -6932.2001

$ws = $wb.Worksheets.Item("WVR")
# WVR row 122
$ws.Range("H122").Value = 6233.2354
$ws.Range("I122").Value = 6264
$ws.Range("J122").Value = 6168.909
$ws.Range("K122").Value = 18792
$ws.Range("L122").Value = 18506.727
$ws.Range("M122").Value = -16342
$ws.Range("N122").Value = -23406.727

# WVR row 126
$ws.Range("H126").Value = 2645.75
$ws.Range("I126").Value = 2645.75
$ws.Range("K126").Value = 7937.25
$ws.Range("M126").Value = -5467.25

# WVR row 132
$ws.Range("H132").Value = 886
$ws.Range("I132").Value = 917.8889
$ws.Range("K132").Value = 2753.6667
$ws.Range("M132").Value = -223.6667000000002
